$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 7: B7, C7, E7 should be numeric values, not text
$ws.Range("B7").Value = 56348
$ws.Range("C7").Value = -1065
$ws.Range("E7").Value = 0

# Add new row 8 with the 2022-01-04 entry (values stored as text, matching
# the original inlineStr formatting style used for incomplete/unfinished rows)
$ws.Range("A8:E8").NumberFormat = "@"
$ws.Range("A8").Value = "2022-01-04"
$ws.Range("B8").Value = "56348.0"
$ws.Range("C8").Value = "-1090.0"
$ws.Range("D8").Value = "-1.93%"
$ws.Range("E8").Value = "0"

# F8/G8 stay blank (same empty "未平倉績效/總績效" columns style as the
# other unfinished rows) - leave untouched rather than writing "", which
# the host normalizes to a cleared cell anyway.

# Reset number format back to default so the new row matches the
# unstyled cells used elsewhere in the sheet (keeps the Text type though)
$ws.Range("A8:E8").Style = "Normal"
